$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage for numeric-looking
# strings (e.g. "0.999", "585.14") so Excel does not silently coerce them
# into real numbers; original cell style/number-format is restored after.
function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$ws.Range("D2").Value = '68.116.97'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '3.336.83'
$ws.Range("E3").Value = '  +0.51%  '
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue $ws.Range("D5") '585.14'
$ws.Range("E5").Value = '  +0.80%  '
Set-TextValue $ws.Range("D6") '177.00'
$ws.Range("E6").Value = '  +1.94%  '
Set-TextValue $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  -0.02%  '
Set-TextValue $ws.Range("D8") '0.593'
$ws.Range("E8").Value = '  +1.44%  '
Set-TextValue $ws.Range("D9") '0.184'
$ws.Range("E9").Value = '  +4.96%  '
Set-TextValue $ws.Range("D10") '0.584'
$ws.Range("E10").Value = '  +1.53%  '
Set-TextValue $ws.Range("D11") '47.98'
$ws.Range("E11").Value = '  +5.69%  '
$ws.Range("E12").Value = '  +2.16%  '
Set-TextValue $ws.Range("D13") '694.66'
$ws.Range("E13").Value = '  +5.04%  '
$ws.Range("D14").Value = '3.880.29'
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = '68.111.82'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").Value = '3.329.43'
$ws.Range("E18").Value = '  +0.54%  '
Set-TextValue $ws.Range("D19") '17.50'
$ws.Range("E19").Value = '  +0.48%  '
Set-TextValue $ws.Range("D20") '11.16'
$ws.Range("E20").Value = '  +2.70%  '
Set-TextValue $ws.Range("D21") '0.895'
$ws.Range("E21").Value = '  +0.96%  '
Set-TextValue $ws.Range("D22") '5.45'
$ws.Range("E22").Value = '  +0.98%  '
Set-TextValue $ws.Range("D23") '16.96'
$ws.Range("E23").Value = '  +0.15%  '
Set-TextValue $ws.Range("D24") '100.69'
$ws.Range("E24").Value = '  +3.56%  '
Set-TextValue $ws.Range("D25") '3.92'
$ws.Range("E25").Value = '  +2.04%  '
$ws.Range("E26").Value = '  +1.00%  '
Set-TextValue $ws.Range("D27") '9.49'
$ws.Range("E27").Value = '  +2.46%  '
Set-TextValue $ws.Range("D28") '33.12'
$ws.Range("E28").Value = '  -0.68%  '
Set-TextValue $ws.Range("D29") '8.53'
$ws.Range("E29").Value = '  +1.82%  '
Set-TextValue $ws.Range("D30") '6.97'
$ws.Range("E30").Value = '  -4.29%  '
Set-TextValue $ws.Range("D31") '567.35'
$ws.Range("E31").Value = '  -2.50%  '
Set-TextValue $ws.Range("D32") '11.05'
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D34").Value = '3.741.09'
$ws.Range("E34").Value = '  +1.00%  '
Set-TextValue $ws.Range("D35") '57.50'
$ws.Range("E35").Value = '  +1.25%  '
$ws.Range("E36").Value = '  -0.04%  '
Set-TextValue $ws.Range("D37") '3.34'
$ws.Range("E37").Value = '  +2.45%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D38") '0.136'
$ws.Range("E38").Value = '  +3.56%  '
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D39") '35.21'
$ws.Range("E39").Value = '  +8.46%  '
Set-TextValue $ws.Range("D40") '3.17'
$ws.Range("E40").Value = '  +2.90%  '
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").Value = '0.0₃0676'
$ws.Range("E42").Value = '  +2.28%  '
Set-TextValue $ws.Range("D43") '0.335'
$ws.Range("E43").Value = '  +0.98%  '
Set-TextValue $ws.Range("D44") '3.26'
$ws.Range("E44").Value = '  +0.56%  '
Set-TextValue $ws.Range("D45") '0.0412'
$ws.Range("E45").Value = '  +1.31%  '
Set-TextValue $ws.Range("D46") '2.66'
$ws.Range("E46").Value = '  +3.35%  '
$ws.Range("E47").Value = '  +1.33%  '
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("E49").Value = '  -0.65%  '
Set-TextValue $ws.Range("D50") '130.90'
$ws.Range("E50").Value = '  +2.79%  '
Set-TextValue $ws.Range("D51") '2.59'
$ws.Range("E51").Value = '  +0.33%  '
